$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content looks like a plain number need the cell
# pre-formatted as Text, otherwise Excel auto-converts the string into a
# float (losing formatting like trailing zeros / thousands separators).

$ws.Range("D2").Value = "43.437.74"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "2.375.71"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.53"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.73"
$ws.Range("E6").Value = "  +3.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.30"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.79"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0814"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.99"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "2.745.43"
$ws.Range("E15").Value = "  +3.36%  "
$ws.Range("E16").Value = "  +5.90%  "
$ws.Range("D17").Value = "2.376.09"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.817"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "43.428.14"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.01"
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0928"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.30"
$ws.Range("E22").Value = "  +3.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.48"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "242.67"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.63"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "26.04"
$ws.Range("E28").Value = "  +8.02%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  +4.81%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.53"
$ws.Range("E30").Value = "  -4.84%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.61"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.28"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.29"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.33"
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +6.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0742"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.67"
$ws.Range("E39").Value = "  +11.15%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.95"
$ws.Range("E40").Value = "  +5.80%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.43"
$ws.Range("E43").Value = "  +4.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.48"
$ws.Range("E44").Value = "  +6.00%  "
$ws.Range("D45").Value = "2.015.52"
$ws.Range("E45").Value = "  +2.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0292"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.20"
$ws.Range("E47").Value = "  +5.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.46"
$ws.Range("E48").Value = "  +6.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "58.11"
$ws.Range("E49").Value = "  +5.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.93"
$ws.Range("E50").Value = "  -3.06%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.580.61"
$ws.Range("E51").Value = "  +2.19%  "
